# Update the "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 4" text to
# "...EMPLOYEE 104" everywhere it appears in the workbook. This text is
# stored as a single shared string referenced from the "header" cell
# (A2 or B2, depending on the sheet layout) on every sheet except "first".
# All occurrences must be updated together so the workbook's shared string
# table collapses back down to the original number of unique entries
# (the stale text is not left behind as an orphan shared string).

$wb = $excel.ActiveWorkbook

$oldText = "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 4"
$newText = "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 104"

# Sheet name -> cell address that holds the banner text.
$targets = @{
    "ResetEmployeeData"            = "A2"
    "UpdteTaxCodeAndAnualSalaryM1" = "A2"
    "ProcessPayrollForApril"       = "B2"
    "ProcessFinalPayrollForApril"  = "B2"
    "TestAprilReports"             = "B2"
    "ProcessPayrollForMay"         = "B2"
    "ProcessFinalPayrollForMay"    = "B2"
    "TestMayReports"               = "B2"
    "UpdteTaxCodeAndAnualSalaryM3" = "A2"
    "ProcessPayrollForJune"        = "B2"
    "ProcessFinalPayrollForJune"   = "B2"
    "TestJuneReports"              = "B2"
    "ProcessPayrollForJuly"        = "B2"
    "ProcessFinalPayrollForJuly"   = "B2"
    "TestJulyReports"              = "B2"
    "ProcessPayrollForAug"         = "B2"
    "ProcessFinalPayrollForAug"    = "B2"
    "TestAugReports"               = "B2"
    "ProcessPayrollForSep"         = "B2"
    "ProcessFinalPayrollForSep"    = "B2"
    "TestSepReports"               = "B2"
    "ProcessPayrollForOct"         = "B2"
    "ProcessFinalPayrollForOct"    = "B2"
    "TestOctReports"               = "B2"
    "ProcessPayrollForNov"         = "B2"
    "ProcessFinalPayrollForNov"    = "B2"
    "TestNovReports"               = "B2"
    "UpdteTaxCodeAndAnualSalaryM9" = "A2"
    "ProcessPayrollForDec"         = "B2"
    "ProcessFinalPayrollForDec"    = "B2"
    "TestDecReports"               = "B2"
    "UpdteTaxCodeAndAnualSalaryM10"= "A2"
    "ProcessPayrollForJan"         = "B2"
    "ProcessFinalPayrollForJan"    = "B2"
    "TestJanReports"               = "B2"
    "ProcessPayrollForFeb"         = "B2"
    "ProcessFinalPayrollForFeb"    = "B2"
    "TestFebReports"               = "B2"
}

foreach ($sheetName in $targets.Keys) {
    $cellAddress = $targets[$sheetName]
    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Range($cellAddress)
    $cell.Value = $newText
}
